$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row cells: "<name>_old" -> "<name>_FV2410",
#        "<name>_new" -> "<name>_FV2504" (the "diff" column is untouched) ---
$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1
$headerRowIndex = $firstRow

for ($c = $firstCol; $c -le $lastCol; $c++) {
  $cell = $ws.Cells.Item($headerRowIndex, $c)
  $text = $cell.Value2
  if ($text -like "*_old") {
    $cell.Value2 = ($text -replace "_old$", "_FV2410")
  } elseif ($text -like "*_new") {
    $cell.Value2 = ($text -replace "_new$", "_FV2504")
  }
}

# --- 2. Turn the data range into an Excel Table (ListObject) ---
# Excel captures any pre-existing ad-hoc header formatting (bold/fill/border)
# as a dxf/headerRowDxfId the first time a table is created over already
# formatted cells. Stash that formatting, clear it, build the table, then
# restore the formatting so the header look (and styles.xml) is unchanged.
$headerRange = $ws.Range($ws.Cells.Item($headerRowIndex, $firstCol), $ws.Cells.Item($headerRowIndex, $lastCol))
$stashRange = $ws.Range($ws.Cells.Item($lastRow + 6, $firstCol), $ws.Cells.Item($lastRow + 6, $lastCol))
$headerRange.Copy()
$stashRange.PasteSpecial(-4122)
$headerRange.ClearFormats()

$dataRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.TableStyle = ""

$stashRange.Copy()
$headerRange.PasteSpecial(-4122)
$stashRange.Clear()

# --- 3. Freeze the header row ---
$ws.Activate()
$ws.Cells.Item($headerRowIndex + 1, $firstCol).Select()
$excel.ActiveWindow.FreezePanes = $true
